$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the existing 3 data rows (rows 3-5) before they get overwritten,
# since the two new tracks are inserted above them (final rows 5-7).
$row3 = @($ws.Range("A3").Value(), $ws.Range("B3").Value(), $ws.Range("C3").Value(), $ws.Range("D3").Value())
$row4 = @($ws.Range("A4").Value(), $ws.Range("B4").Value(), $ws.Range("C4").Value(), $ws.Range("D4").Value())
$row5 = @($ws.Range("A5").Value(), $ws.Range("B5").Value(), $ws.Range("C5").Value(), $ws.Range("D5").Value())

# Push the existing rows down two slots: 3->5, 4->6, 5->7.
$ws.Range("A7").Value = $row5[0]
$ws.Range("B7").Value = $row5[1]
$ws.Range("C7").Value = $row5[2]
$ws.Range("D7").Value = $row5[3]

$ws.Range("A6").Value = $row4[0]
$ws.Range("B6").Value = $row4[1]
$ws.Range("C6").Value = $row4[2]
$ws.Range("D6").Value = $row4[3]

$ws.Range("A5").Value = $row3[0]
$ws.Range("B5").Value = $row3[1]
$ws.Range("C5").Value = $row3[2]
$ws.Range("D5").Value = $row3[3]

# Match the text-style formatting used by the "Name" column so the new
# rows carry the same cell style (s="1") as the rest of the table.
$ws.Range("A5:A7").NumberFormat = "@"

# New row 1: "Exceed the Sky"
$ws.Range("A3").Value = "Exceed the Sky"
$ws.Range("A3").NumberFormat = "@"
$ws.Range("B3").Value = 21943
$ws.Range("C3").Value = 21943
$ws.Range("D3").Value = 61314

# New row 2: "Junk Stereo Concept"
$ws.Range("A4").Value = "Junk Stereo Concept"
$ws.Range("A4").NumberFormat = "@"
$ws.Range("B4").Value = 20
$ws.Range("C4").Value = 15360
$ws.Range("D4").Value = 61440

# Grow the XML-mapped table (and its autofilter) to cover the two new rows.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A2:D7"))
